$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.442.64'
$ws.Range('E2').Value = '  -0.88%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.054.95'
$ws.Range('E3').Value = '  -1.29%  '

$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.00'
$ws.Range('E5').Value = '  -1.55%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('E7').Value = '  +0.12%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.97'
$ws.Range('E8').Value = '  -3.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  -2.34%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0769'
$ws.Range('E10').Value = '  -2.26%  '

$ws.Range('E11').Value = '  +1.51%  '

$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.371.86'
$ws.Range('E12').Value = '  -0.75%  '

$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.73'
$ws.Range('E13').Value = '  -0.95%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.56'
$ws.Range('E14').Value = '  -2.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.754'
$ws.Range('E15').Value = '  -2.77%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.24'
$ws.Range('E16').Value = '  -2.26%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.079.64'
$ws.Range('E17').Value = '  -0.45%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.436.86'
$ws.Range('E18').Value = '  -0.78%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.90'
$ws.Range('E19').Value = '  -2.46%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.88'
$ws.Range('E20').Value = '  -4.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0819'
$ws.Range('E21').Value = '  -2.10%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.94'
$ws.Range('E22').Value = '  -1.13%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.35'
$ws.Range('E24').Value = '  +3.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.33'
$ws.Range('E25').Value = '  -3.50%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.43'
$ws.Range('E26').Value = '  +3.05%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.18'
$ws.Range('E27').Value = '  -1.70%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.130'
$ws.Range('E28').Value = '  -4.00%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.14'
$ws.Range('E29').Value = '  -1.94%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.35'
$ws.Range('E30').Value = '  -4.39%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.121'
$ws.Range('E31').Value = '  +0.15%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  -3.84%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0622'
$ws.Range('E33').Value = '  -2.30%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.54'
$ws.Range('E34').Value = '  -3.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.44'
$ws.Range('E35').Value = '  -1.85%  '

$ws.Range('E36').Value = '  +0.06%  '

$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.20%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.25'
$ws.Range('E38').Value = '  -4.41%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.23'
$ws.Range('E39').Value = '  -2.36%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0226'
$ws.Range('E40').Value = '  +4.79%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.83'
$ws.Range('E41').Value = '  -1.92%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0950'
$ws.Range('E42').Value = '  -2.98%  '

$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.89'
$ws.Range('E43').Value = '  +0.11%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.472.86'
$ws.Range('E44').Value = '  +2.13%  '

$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.19'
$ws.Range('E45').Value = '  +2.87%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.34'
$ws.Range('E46').Value = '  -2.10%  '

$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.06'
$ws.Range('E47').Value = '  -2.94%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.03'
$ws.Range('E48').Value = '  -3.59%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.15'
$ws.Range('E49').Value = '  -3.86%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('E50').Value = '  -2.44%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.255.18'
$ws.Range('E51').Value = '  -0.76%  '
